$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated financial figures (numeric cell value changes)
$updates = @{
    "D8" = 17800
    "H8" = 196000
    "I8" = 205400
    "J8" = 175100
    "D9" = 5700
    "G9" = 40700
    "H9" = 138400
    "I9" = 136900
    "J9" = 115000
    "D10" = 12100
    "G10" = -40700
    "H10" = 57600
    "I10" = 68500
    "J10" = 60100
    "D17" = 10700
    "G17" = 700
    "H17" = 175800
    "I17" = 170400
    "J17" = 125500
    "D18" = 7000
    "G18" = -700
    "H18" = 20200
    "I18" = 35000
    "J18" = 49600
    "D21" = 7200
    "H21" = 21300
    "I21" = 35700
    "J21" = 49800
    "D23" = 7200
    "G23" = -700
    "H23" = 20200
    "I23" = 35000
    "J23" = 49600
    "D24" = 2200
    "H24" = 5800
    "I24" = 8900
    "J24" = 12200
    "D26" = 5000
    "G26" = -700
    "H26" = 14400
    "I26" = 26000
    "J26" = 37400
    "D27" = 4000
    "G27" = -700
    "H27" = 14400
    "I27" = 26000
    "J27" = 37400
    "E29" = -50700
    "F29" = -88800
    "G29" = -24700
    "D33" = 4000
    "E33" = -51200
    "F33" = -89100
    "G33" = -25300
    "H33" = 14400
    "I33" = 26000
    "J33" = 37400
    "D35" = 4000
    "E35" = -51200
    "F35" = -89100
    "G35" = -25300
    "H35" = 14400
    "I35" = 26000
    "J35" = 37400
    "F41" = 30200
    "G41" = 146200
    "H41" = 119700
    "I41" = 162700
    "J41" = 153100
    "G42" = 10400
    "H42" = 19300
    "I42" = 7400
    "D43" = 123300
    "F43" = 96400
    "G43" = 49900
    "H43" = 88500
    "I43" = 54000
    "J43" = 51800
    "G44" = 14500
    "H44" = 4800
    "D45" = 2300
    "F45" = 6800
    "G45" = 7100
    "H45" = 12100
    "I45" = 13500
    "J45" = 12800
    "D46" = 128800
    "F46" = 135600
    "G46" = 228200
    "H46" = 244400
    "I46" = 239100
    "J46" = 217300
    "D48" = 100
    "F49" = 2600
    "H49" = 1300
    "I49" = 1900
    "E52" = 81000
    "G52" = 9400
    "D54" = 129000
    "E54" = 81000
    "F54" = 138700
    "G54" = 239400
    "H54" = 248400
    "I54" = 244600
    "J54" = 221300
    "G57" = 6700
    "H57" = 3800
    "I57" = 13400
    "J57" = 12400
    "D58" = 33600
    "D59" = 5700
    "E59" = 23400
    "F59" = 29100
    "G59" = 35100
    "H59" = 21900
    "I59" = 23000
    "J59" = 31400
    "D60" = 39300
    "E60" = 23400
    "F60" = 30100
    "G60" = 41800
    "H60" = 25700
    "I60" = 36400
    "J60" = 38200
    "D66" = 57300
    "E66" = 23400
    "F66" = 30100
    "G66" = 41800
    "H66" = 25700
    "I66" = 36400
    "J66" = 38200
    "D72" = 14000
    "E72" = -16900
    "F72" = 34300
    "G72" = 123400
    "H72" = 148600
    "I72" = 134000
    "J72" = 107700
    "D76" = 71700
    "E76" = 57600
    "F76" = 108600
    "G76" = 197500
    "H76" = 222600
    "I76" = 208200
    "J76" = 183100
    "D81" = 4000
    "E81" = -51200
    "F81" = -89100
    "G81" = -25300
    "H81" = 14400
    "I81" = 26000
    "J81" = 37400
    "H83" = 1100
    "D89" = 6600
    "E89" = -5100
    "F89" = -107300
    "G89" = 16800
    "H89" = -33200
    "I89" = 18700
    "J89" = 30900
    "D94" = -20200
    "E94" = -19900
    "F94" = -8700
    "G94" = 9800
    "H94" = -9800
    "I94" = -8000
    "D100" = 2500
    "J100" = -1600
    "D102" = -11100
    "E102" = -25000
    "F102" = -116000
    "G102" = 26500
    "H102" = -43000
    "I102" = 9500
    "J102" = 25100
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Cells that switch from numeric 0 to the "NA" text marker
$naCells = @("G91", "H91", "I91", "J91")
foreach ($addr in $naCells) {
    $ws.Range($addr).Value = "NA"
}

